$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.578.40"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "'3.477.76"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'589.08"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").Value = "'168.14"
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("D7").Value = "'0.606"
$ws.Range("E7").Value = "  -1.76%  "
$ws.Range("D8").Value = "'3.472.87"
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").Value = "'6.79"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").Value = "'0.572"
$ws.Range("E12").Value = "  -4.11%  "
$ws.Range("D13").Value = "'46.52"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").Value = "'4.032.19"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "'614.15"
$ws.Range("E16").Value = "  -10.15%  "
$ws.Range("D17").Value = "'8.30"
$ws.Range("D18").Value = "'3.472.98"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'68.595.13"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").Value = "'11.10"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'0.870"
$ws.Range("E23").Value = "  -4.07%  "
$ws.Range("E24").Value = "  -4.40%  "
$ws.Range("D25").Value = "'95.62"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").Value = "'3.77"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  -1.86%  "
$ws.Range("D29").Value = "'9.09"
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("D30").Value = "'32.76"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").Value = "'8.39"
$ws.Range("E31").Value = "  -4.78%  "
$ws.Range("D32").Value = "'3.08"
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").Value = "'1.32"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("E34").Value = "  -6.08%  "
$ws.Range("D35").Value = "'573.50"
$ws.Range("E35").Value = "  -0.92%  "
$ws.Range("D36").Value = "'10.68"
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").Value = "'3.49"
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("D38").Value = "'56.91"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  -3.57%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "'0.137"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'0.0437"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'3.387.15"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("E44").Value = "  -4.20%  "
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("D47").Value = "'2.82"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("E49").Value = "  -3.41%  "
$ws.Range("D50").Value = "'132.41"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").Value = "'5.55"
$ws.Range("E51").Value = "  +9.32%  "
